$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at D: "Battery total Capacity (MWh)"
$ws.Columns.Item(4).Insert()
$ws.Range("D1").Value = "Battery total Capacity (MWh)"

# Update row 2 values to reflect the new 5-year data run
$ws.Range("A2").Value = 446.5698322765596
$ws.Range("B2").Value = 0
$ws.Range("C2").Value = 250
$ws.Range("D2").Value = 2500
$ws.Range("E2").Value = 4459.25831949529
$ws.Range("F2").Value = 5459.25831949529
$ws.Range("G2").Value = 21972995369.31326
$ws.Range("H2").Value = 45.00000000000043
$ws.Range("I2").Value = 4927500.000000047
$ws.Range("J2").Value = [double]"3.364582101824124e-12"
$ws.Range("K2").Value = 5041414.416467906
$ws.Range("L2").Value = 10950000
$ws.Range("M2").Value = 1000
$ws.Range("N2").Value = 21587498190.05507
